# Apply the edits described by the commit "added new elements for mailing address":
#  1. Update the letter date from September 19 to September 21, 2025.
#  2. Split the mailing address onto two lines/paragraphs:
#       "909 Story Road, San Jose CA 95122"
#     becomes
#       "909 Story Road"
#       "San Jose, CA 95122"
#  3. Remove the now-superfluous blank "NoSpacing" paragraph that sat
#     right after the "...Board of Directors" paragraph.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. Date update ---------------------------------------------------
$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($f.Start, $f.End)
$body = '<w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">September 21, 2025</w:t></w:r></w:p></w:body>'
$target.InsertXML($pkgHeader + $body + $pkgFooter) | Out-Null

# --- 2. Split the address line into two paragraphs --------------------
$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Execute("909 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($f.Start, $f.End)
$body = '<w:body>' +
        '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">909 Story Road</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">San Jose, CA 95122</w:t></w:r></w:p>' +
        '</w:body>'
$target.InsertXML($pkgHeader + $body + $pkgFooter) | Out-Null

# --- 3. Remove the blank NoSpacing paragraph after "Board of Directors" --
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*Board of Directors`r") {
        $next = $p.Next()
        if ($next -ne $null) {
            $next.Range.Delete() | Out-Null
        }
        break
    }
}
